$d = $word.ActiveDocument

# Locate the unique "12/9 at 1:30pm" table-cell text and replace it so that
# the time portion ("12/9 at 1:10") and the trailing "pm" end up as two
# separate (but identically formatted) runs, matching how Word splits a run
# when only part of its text is retyped.

$target = $d.Content.Duplicate
$found = $target.Find.Execute("12/9 at 1:30pm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '12/9 at 1:30pm' in the document"
}

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p w:rsidR="006F38B7" w:rsidRPr="003344AD" w:rsidRDefault="005E78DB" w:rsidP="003344AD">
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:left="0"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>12/9 at 1:10</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>pm</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.InsertXML($xml)
